$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.409.31'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '3.495.99'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').Value = '  -0.08%  '
$cell = $ws.Range('D5')
$cell.Formula = '="586.55"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E5').Value = '  +0.01%  '
$cell = $ws.Range('D6')
$cell.Formula = '="134.87"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E6').Value = '  +1.73%  '
$ws.Range('D7').Value = '3.494.90'
$ws.Range('E7').Value = '  -0.43%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('E10').Value = '  -0.39%  '
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('E12').Value = '  -3.31%  '
$ws.Range('D13').Value = '4.088.52'
$ws.Range('E14').Value = '  +0.37%  '
$ws.Range('E15').Value = '  +1.13%  '
$ws.Range('D16').Value = '3.493.58'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '64.417.67'
$ws.Range('E17').Value = '  -0.39%  '
$cell = $ws.Range('D18')
$cell.Formula = '="25.12"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E18').Value = '  -9.69%  '
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('E20').Value = '  -0.45%  '
$cell = $ws.Range('D21')
$cell.Formula = '="13.75"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E21').Value = '  -4.65%  '
$cell = $ws.Range('D22')
$cell.Formula = '="385.30"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E22').Value = '  -1.34%  '
$cell = $ws.Range('D23')
$cell.Formula = '="0.568"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E23').Value = '  -1.85%  '
$ws.Range('D24').Value = '3.634.76'
$ws.Range('E24').Value = '  -0.56%  '
$cell = $ws.Range('D25')
$cell.Formula = '="74.09"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E25').Value = '  +0.58%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E27').Value = '  +3.54%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range('D28')
$cell.Formula = '="7.46"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell = $ws.Range('D29')
$cell.Formula = '="1.54"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E29').Value = '  -2.61%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('E31').Value = '  -1.59%  '
$cell = $ws.Range('D32')
$cell.Formula = '="8.21"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('D33').Value = '3.514.00'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('E35').Value = '  +0.48%  '
$ws.Range('E36').Value = '  -2.30%  '
$ws.Range('E37').Value = '  +0.55%  '
$cell = $ws.Range('D38')
$cell.Formula = '="1.54"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E38').Value = '  -3.56%  '
$cell = $ws.Range('D39')
$cell.Formula = '="6.84"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E39').Value = '  -2.22%  '
$cell = $ws.Range('D40')
$cell.Formula = '="162.57"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E40').Value = '  -5.04%  '
$ws.Range('E41').Value = '  -3.34%  '
$cell = $ws.Range('D42')
$cell.Formula = '="0.805"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E42').Value = '  -1.06%  '
$cell = $ws.Range('D43')
$cell.Formula = '="25.83"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E43').Value = '  -2.90%  '
$ws.Range('E45').Value = '  -0.57%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range('D46')
$cell.Formula = '="4.41"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E46').Value = '  +0.44%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$cell = $ws.Range('D47')
$cell.Formula = '="1.21"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E47').Value = '  -0.22%  '
$cell = $ws.Range('D48')
$cell.Formula = '="1.65"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('D49').Value = '2.474.35'
$ws.Range('E49').Value = '  +1.05%  '
$ws.Range('E50').Value = '  -1.86%  '
$cell = $ws.Range('D51')
$cell.Formula = '="0.906"'
$cell.Copy()
$cell.PasteSpecial("xlPasteValues")
$ws.Range('E51').Value = '  +0.84%  '
